$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete the data row for VICTOR JOSE LOPEZ GALVAN (row 16), shifting
# LIZETH MARTINEZ DIAZ's row (17) and the signature rows (22,23) up by one.
$ws.Rows("16").Delete()

# Update "VALOR MORA" total (E11) now that only one worker remains.
$ws.Range("E11").Value = 16000

# Update "Cant. Trabajadores" (C13) count.
$ws.Range("C13").Value = 1

# Narrow column D now that the longer name has been removed (best-fit
# re-calculation after VICTOR JOSE LOPEZ GALVAN's longer name is gone).
$ws.Columns("D").ColumnWidth = 20.5

$wb.Save()
